$wb = $excel.ActiveWorkbook

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item('Overall')
$ws.Cells.Item(1, 1).NumberFormat = '@'
$ws.Cells.Item(1, 1).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 2).NumberFormat = '@'
$ws.Cells.Item(1, 2).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 3).NumberFormat = '@'
$ws.Cells.Item(1, 3).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 4).NumberFormat = '@'
$ws.Cells.Item(1, 4).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 5).NumberFormat = '@'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = '73.54%'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '514'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '$1,256,309,593'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '9.93%'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '-26.82%'

# ---- Sheet: County ----
$ws = $wb.Worksheets.Item('County')
$ws.Cells.Item(1, 1).NumberFormat = '@'
$ws.Cells.Item(1, 1).Value = 'Geography'
$ws.Cells.Item(1, 2).NumberFormat = '@'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).NumberFormat = '@'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).NumberFormat = '@'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).NumberFormat = '@'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).NumberFormat = '@'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = 'United States'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '67.35%'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '103,475'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$267,700,640,005'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '9.05%'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-12.83%'
$ws.Cells.Item(3, 1).NumberFormat = '@'
$ws.Cells.Item(3, 1).Value = 'Delaware'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '73.54%'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '514'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$1,256,309,593'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '9.93%'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-26.82%'
$ws.Cells.Item(4, 1).NumberFormat = '@'
$ws.Cells.Item(4, 1).Value = 'Kent County'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '83.75%'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '80'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$254,877,640'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '7.32%'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '-47.74%'
$ws.Cells.Item(5, 1).NumberFormat = '@'
$ws.Cells.Item(5, 1).Value = 'New Castle County'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '71.26%'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '334'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$896,432,426'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '8.76%'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-22.78%'
$ws.Cells.Item(6, 1).NumberFormat = '@'
$ws.Cells.Item(6, 1).Value = 'Sussex County'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '73.00%'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '100'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$104,999,527'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '17.21%'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '-29.46%'

# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item('Congressional District')
$ws.Cells.Item(1, 1).NumberFormat = '@'
$ws.Cells.Item(1, 1).Value = 'Geography'
$ws.Cells.Item(1, 2).NumberFormat = '@'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).NumberFormat = '@'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).NumberFormat = '@'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).NumberFormat = '@'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).NumberFormat = '@'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = 'United States'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '67.35%'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '103,475'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$267,700,640,005'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '9.05%'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-12.83%'
$ws.Cells.Item(3, 1).NumberFormat = '@'
$ws.Cells.Item(3, 1).Value = 'Delaware'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '73.54%'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '514'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$1,256,309,593'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '9.93%'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-26.82%'
$ws.Cells.Item(4, 1).NumberFormat = '@'
$ws.Cells.Item(4, 1).Value = 'Congressional District (at Large)'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '73.54%'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '514'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$1,256,309,593'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '9.93%'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '-26.82%'

# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item('Size')
$ws.Cells.Item(1, 1).NumberFormat = '@'
$ws.Cells.Item(1, 1).Value = 'Size'
$ws.Cells.Item(1, 2).NumberFormat = '@'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).NumberFormat = '@'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).NumberFormat = '@'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).NumberFormat = '@'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).NumberFormat = '@'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = 'Between $100K and $499K'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '75.00%'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '156'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$28,746,625'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '7.66%'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-35.66%'
$ws.Cells.Item(3, 1).NumberFormat = '@'
$ws.Cells.Item(3, 1).Value = 'Between $1M and $4.99M'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '70.86%'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '151'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$159,238,036'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '14.27%'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-31.12%'
$ws.Cells.Item(4, 1).NumberFormat = '@'
$ws.Cells.Item(4, 1).Value = 'Between $500K and $999K'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '76.92%'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '91'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$39,293,365'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '13.41%'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '-30.20%'
$ws.Cells.Item(5, 1).NumberFormat = '@'
$ws.Cells.Item(5, 1).Value = 'Between $5M and $9.99M'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '84.62%'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '39'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$148,811,089'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '5.63%'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-36.23%'
$ws.Cells.Item(6, 1).NumberFormat = '@'
$ws.Cells.Item(6, 1).Value = 'Greater than $10M'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '69.39%'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '49'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$878,068,273'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '5.16%'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '-12.65%'
$ws.Cells.Item(7, 1).NumberFormat = '@'
$ws.Cells.Item(7, 1).Value = 'Less than $100K'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '60.71%'
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = '28'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '$2,152,205'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '11.26%'
$ws.Cells.Item(7, 6).NumberFormat = '@'
$ws.Cells.Item(7, 6).Value = '-23.55%'
$ws.Cells.Item(8, 1).NumberFormat = '@'
$ws.Cells.Item(8, 1).Value = 'Total'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '73.54%'
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = '514'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '$1,256,309,593'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '9.93%'
$ws.Cells.Item(8, 6).NumberFormat = '@'
$ws.Cells.Item(8, 6).Value = '-26.82%'

# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item('Subsector')
$ws.Cells.Item(1, 1).NumberFormat = '@'
$ws.Cells.Item(1, 1).Value = 'Subsector'
$ws.Cells.Item(1, 2).NumberFormat = '@'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).NumberFormat = '@'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).NumberFormat = '@'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).NumberFormat = '@'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).NumberFormat = '@'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = 'Arts, Culture, and Humanities'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '68.75%'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '48'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$32,126,193'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '13.50%'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-10.17%'
$ws.Cells.Item(3, 1).NumberFormat = '@'
$ws.Cells.Item(3, 1).Value = 'Education (Excluding Universities)'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '72.13%'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '61'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$197,388,696'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '5.99%'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-56.44%'
$ws.Cells.Item(4, 1).NumberFormat = '@'
$ws.Cells.Item(4, 1).Value = 'Environment and Animals'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '43.48%'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '23'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$249,975,545'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '25.58%'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '1.20%'
$ws.Cells.Item(5, 1).NumberFormat = '@'
$ws.Cells.Item(5, 1).Value = 'Health (Excluding Hospitals)'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '81.82%'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '33'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$44,106,600'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '1.35%'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-43.98%'
$ws.Cells.Item(6, 1).NumberFormat = '@'
$ws.Cells.Item(6, 1).Value = 'Hospitals'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '33.33%'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '3'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$37,023,317'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '9.32%'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '5.61%'
$ws.Cells.Item(7, 1).NumberFormat = '@'
$ws.Cells.Item(7, 1).Value = 'Human Services'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '76.76%'
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = '185'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '$181,205,892'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '13.42%'
$ws.Cells.Item(7, 6).NumberFormat = '@'
$ws.Cells.Item(7, 6).Value = '-39.73%'
$ws.Cells.Item(8, 1).NumberFormat = '@'
$ws.Cells.Item(8, 1).Value = 'International, Foreign Affairs'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '75.00%'
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = '4'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '$1,840,864'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '0.50%'
$ws.Cells.Item(8, 6).NumberFormat = '@'
$ws.Cells.Item(8, 6).Value = '-60.23%'
$ws.Cells.Item(9, 1).NumberFormat = '@'
$ws.Cells.Item(9, 1).Value = 'Public, Societal Benefit'
$ws.Cells.Item(9, 2).NumberFormat = '@'
$ws.Cells.Item(9, 2).Value = '70.59%'
$ws.Cells.Item(9, 3).NumberFormat = '@'
$ws.Cells.Item(9, 3).Value = '34'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '$33,197,123'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '11.36%'
$ws.Cells.Item(9, 6).NumberFormat = '@'
$ws.Cells.Item(9, 6).Value = '-15.57%'
$ws.Cells.Item(10, 1).NumberFormat = '@'
$ws.Cells.Item(10, 1).Value = 'Religion Related'
$ws.Cells.Item(10, 2).NumberFormat = '@'
$ws.Cells.Item(10, 2).Value = '80.00%'
$ws.Cells.Item(10, 3).NumberFormat = '@'
$ws.Cells.Item(10, 3).Value = '5'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '$751,535'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '16.97%'
$ws.Cells.Item(10, 6).NumberFormat = '@'
$ws.Cells.Item(10, 6).Value = '-9.69%'
$ws.Cells.Item(11, 1).NumberFormat = '@'
$ws.Cells.Item(11, 1).Value = 'Unclassified'
$ws.Cells.Item(11, 2).NumberFormat = '@'
$ws.Cells.Item(11, 2).Value = '77.19%'
$ws.Cells.Item(11, 3).NumberFormat = '@'
$ws.Cells.Item(11, 3).Value = '114'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '$192,545,718'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '8.92%'
$ws.Cells.Item(11, 6).NumberFormat = '@'
$ws.Cells.Item(11, 6).Value = '-19.66%'
$ws.Cells.Item(12, 1).NumberFormat = '@'
$ws.Cells.Item(12, 1).Value = 'Universities'
$ws.Cells.Item(12, 2).NumberFormat = '@'
$ws.Cells.Item(12, 2).Value = '50.00%'
$ws.Cells.Item(12, 3).NumberFormat = '@'
$ws.Cells.Item(12, 3).Value = '4'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '$286,148,110'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '12.51%'
$ws.Cells.Item(12, 6).NumberFormat = '@'
$ws.Cells.Item(12, 6).Value = '1.63%'
$ws.Cells.Item(13, 1).NumberFormat = '@'
$ws.Cells.Item(13, 1).Value = 'Total'
$ws.Cells.Item(13, 2).NumberFormat = '@'
$ws.Cells.Item(13, 2).Value = '73.54%'
$ws.Cells.Item(13, 3).NumberFormat = '@'
$ws.Cells.Item(13, 3).Value = '514'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '$1,256,309,593'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '9.93%'
$ws.Cells.Item(13, 6).NumberFormat = '@'
$ws.Cells.Item(13, 6).Value = '-26.82%'
